# Weekly update: insert a new price record as row 32 (new fecha/week),
# pushing the existing rows 32-77 down to 33-78.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("32:32").Insert()

$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44671
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108004
$ws.Range("J32").Value = "Papaya"
$ws.Range("K32").Value = "Cultivar IV Región"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 30
$ws.Range("N32").Value = 25000
$ws.Range("O32").Value = 25000
$ws.Range("P32").Value = 25000
$ws.Range("Q32").Value = "$/bandeja 10 kilos"
$ws.Range("R32").Value = "Provincia del Elquí"
$ws.Range("S32").Value = 2500
$ws.Range("T32").Value = 10
